$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-17: Text, X, Y, Width, Height
$data = @(
    @("February ", 908, 644, 83, 23),
    @("Revolution ", 152, 669, 96, 23),
    @(". ", 248, 669, 11, 23),
    @("The ", 259, 669, 39, 23),
    @("second ", 298, 669, 68, 23),
    @("was ", 366, 669, 40, 23),
    @("the ", 406, 669, 32, 23),
    @("October ", 438, 669, 74, 23),
    @("Revolution.", 512, 669, 96, 23),
    @("czarist ", 501, 710.4, 61, 23),
    @("government.", 562, 710.4, 107, 23),
    @("Russian ", 819, 1324.4, 75, 23),
    @("Civil ", 894, 1324.4, 42, 23),
    @("War. ", 152, 1349.4, 45, 23),
    @("Soviet ", 491, 1349.4, 59, 23),
    @("Union ", 550, 1349.4, 56, 23)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
